$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CONDITION row (B8): the DRL snippet template is split so the comparison
# operator/literal moves out into its own shared column, leaving just the
# field name here.
$ws.Range("B8").Value = "insuredId"

# Rows 10/11 (rule1/rule2 parameter column): the raw id values are now
# quoted string literals instead of bare identifiers.
$ws.Range("B10").Value = '"A223456789"'
$ws.Range("B11").Value = '"A223456123"'

# Move/restore the cursor selection to B12, matching the saved workbook view.
$ws.Range("B12").Select()
